$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the old "TCP" row (row 10). This shifts old rows 11-18 up to
#    10-17, which carries along their correct row-height / thick border
#    formatting for free.
# ---------------------------------------------------------------------------
$ws.Rows(10).Delete()

# ---------------------------------------------------------------------------
# 2) Insert two new columns before column C for lower_bound / upper_bound.
#    Old column C (average_best) becomes E, old column D (price) becomes F.
# ---------------------------------------------------------------------------
$ws.Range("C1:D1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "lower_bound"
$ws.Range("D1").Value = "upper_bound"

# ---------------------------------------------------------------------------
# 4) Fill in lower_bound / upper_bound for the treatment rows (2-9)
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 0.064
$ws.Range("D2").Value = 0.31

$ws.Range("C3").Value = 0.1
$ws.Range("D3").Value = 0.56

$ws.Range("C4").Value = 0.01202
$ws.Range("D4").Value = 0.1424

$ws.Range("C5").Value = 0.0144
$ws.Range("D5").Value = 0.0778

$ws.Range("C6").Value = 0.016
$ws.Range("D6").Value = 0.162

$ws.Range("C7").Value = 0.00212
$ws.Range("D7").Value = 0.0286

$ws.Range("C8").Value = 0.0262
$ws.Range("D8").Value = 0.0906

$ws.Range("C9").Value = 0.078
$ws.Range("D9").Value = 0.266

# Re-apply the numeric-cell style (style of E column, e.g. E2) to the new
# C:D cells on rows 2-9, since column-insert copies the left neighbour's
# (name column) formatting by default.
for ($r = 2; $r -le 9; $r++) {
    $ws.Range("E$r").Copy()
    $ws.Range("C$r`:D$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) Fill in the three repeated-bound columns (C, D, E) for rows 10-16
#    (laser_mild .. surgery_sev) using the existing E value (shifted from
#    the old column C) as the source.
# ---------------------------------------------------------------------------
for ($r = 10; $r -le 16; $r++) {
    $v = $ws.Range("E$r").Value
    $ws.Range("C$r").Value = $v
    $ws.Range("D$r").Value = $v
    $ws.Range("E$r").Copy()
    $ws.Range("C$r`:D$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6) Row 17 (surgery_vi) needs to lose its inherited thick border / custom
#    height (the source template file always renders the very last data
#    row with plain/default height), so recreate it from scratch.
# ---------------------------------------------------------------------------
$lastValue = $ws.Range("A17").Value
$ws.Rows(17).Delete()
$ws.Rows(17).Insert()
$ws.Range("A17").Value = $lastValue
$ws.Range("C17").Value = 0.22000000000000003
$ws.Range("D17").Value = 0.22000000000000003
$ws.Range("E17").Value = 0.22000000000000003

$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A16").Copy()
$ws.Range("C17`:E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 7) View state: selection
# ---------------------------------------------------------------------------
$ws.Range("C20").Select()
